$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: BGC0000669
$ws.Cells.Item(2, 1).Value = 0
$ws.Cells.Item(2, 2).Value = 'BGC0000669'
$ws.Cells.Item(2, 3).Value = 'None'
$ws.Cells.Item(2, 4).Value = 'None'
$ws.Cells.Item(2, 5).Value = 'Arabidopsis thaliana'
$ws.Cells.Item(2, 6).Value = 17023645
$ws.Cells.Item(2, 7).Value = 17058245
$ws.Cells.Item(2, 8).Value = 'None'
$ws.Cells.Item(2, 9).Value = 5
$ws.Cells.Item(2, 10).Value = 'Terpene'
$ws.Cells.Item(2, 11).Value = 'Arabidopsis thaliana chromosome 5 sequence'
$ws.Cells.Item(2, 12).Value = 'marneral'

# Row 3: BGC0000670
$ws.Cells.Item(3, 1).Value = 1
$ws.Cells.Item(3, 2).Value = 'BGC0000670'
$ws.Cells.Item(3, 3).Value = 'None'
$ws.Cells.Item(3, 4).Value = 'None'
$ws.Cells.Item(3, 5).Value = 'Arabidopsis thaliana'
$ws.Cells.Item(3, 6).Value = 19428887
$ws.Cells.Item(3, 7).Value = 19461689
$ws.Cells.Item(3, 8).Value = 'None'
$ws.Cells.Item(3, 9).Value = 4
$ws.Cells.Item(3, 10).Value = 'Terpene'
$ws.Cells.Item(3, 11).Value = 'Arabidopsis thaliana chromosome 5 sequence'
$ws.Cells.Item(3, 12).Value = 'thaliandiol,'

# Row 4: BGC0000671
$ws.Cells.Item(4, 1).Value = 2
$ws.Cells.Item(4, 2).Value = 'BGC0000671'
$ws.Cells.Item(4, 3).Value = 'None'
$ws.Cells.Item(4, 4).Value = 'None'
$ws.Cells.Item(4, 5).Value = 'Oryza sativa Japonica Group'
$ws.Cells.Item(4, 6).Value = 5310456
$ws.Cells.Item(4, 7).Value = 5479082
$ws.Cells.Item(4, 8).Value = 'None'
$ws.Cells.Item(4, 9).Value = 7
$ws.Cells.Item(4, 10).Value = 'Terpene'
$ws.Cells.Item(4, 11).Value = 'Oryza sativa Japonica Group DNA, chromosome 4, complete sequence, cultivar: Nipponbare'
$ws.Cells.Item(4, 12).Value = 'momilactone'

# Row 5: BGC0000672
$ws.Cells.Item(5, 1).Value = 3
$ws.Cells.Item(5, 2).Value = 'BGC0000672'
$ws.Cells.Item(5, 3).Value = 'None'
$ws.Cells.Item(5, 4).Value = 'None'
$ws.Cells.Item(5, 5).Value = 'Oryza sativa Japonica Group'
$ws.Cells.Item(5, 6).Value = 22520468
$ws.Cells.Item(5, 7).Value = 22764099
$ws.Cells.Item(5, 8).Value = 'None'
$ws.Cells.Item(5, 9).Value = 10
$ws.Cells.Item(5, 10).Value = 'Terpene'
$ws.Cells.Item(5, 11).Value = 'Oryza sativa Japonica Group DNA, chromosome 2, complete sequence, cultivar: Nipponbare'
$ws.Cells.Item(5, 12).Value = 'oryzalides,'

# Row 6: BGC0000810
$ws.Cells.Item(6, 1).Value = 4
$ws.Cells.Item(6, 2).Value = 'BGC0000810'
$ws.Cells.Item(6, 3).Value = 'None'
$ws.Cells.Item(6, 4).Value = 'None'
$ws.Cells.Item(6, 5).Value = 'Zea mays'
$ws.Cells.Item(6, 6).Value = 3003161
$ws.Cells.Item(6, 7).Value = 3267368
$ws.Cells.Item(6, 8).Value = 'None'
$ws.Cells.Item(6, 9).Value = 8
$ws.Cells.Item(6, 10).Value = 'Alkaloid'
$ws.Cells.Item(6, 11).Value = 'Zea mays cultivar B73 chromosome 4'
$ws.Cells.Item(6, 12).Value = 'benzoxazinone'

# Row 7: BGC0001313
$ws.Cells.Item(7, 1).Value = 5
$ws.Cells.Item(7, 2).Value = 'BGC0001313'
$ws.Cells.Item(7, 3).Value = 'None'
$ws.Cells.Item(7, 4).Value = 'None'
$ws.Cells.Item(7, 5).Value = 'Arabidopsis thaliana'
$ws.Cells.Item(7, 6).Value = 8729999
$ws.Cells.Item(7, 7).Value = 8820000
$ws.Cells.Item(7, 8).Value = 'None'
$ws.Cells.Item(7, 9).Value = 18
$ws.Cells.Item(7, 10).Value = 'Terpene'
$ws.Cells.Item(7, 11).Value = 'Arabidopsis thaliana chromosome 4 sequence'
$ws.Cells.Item(7, 12).Value = 'arabidiol-baruol'

# Row 8: BGC0001314
$ws.Cells.Item(8, 1).Value = 6
$ws.Cells.Item(8, 2).Value = 'BGC0001314'
$ws.Cells.Item(8, 3).Value = 'None'
$ws.Cells.Item(8, 4).Value = 'None'
$ws.Cells.Item(8, 5).Value = 'Arabidopsis thaliana'
$ws.Cells.Item(8, 6).Value = 14189999
$ws.Cells.Item(8, 7).Value = 14250000
$ws.Cells.Item(8, 8).Value = 'None'
$ws.Cells.Item(8, 9).Value = 14
$ws.Cells.Item(8, 10).Value = 'Terpene'
$ws.Cells.Item(8, 11).Value = 'Arabidopsis thaliana chromosome 5 sequence'
$ws.Cells.Item(8, 12).Value = 'tirucalla'

# Row 9: BGC0001315
$ws.Cells.Item(9, 1).Value = 7
$ws.Cells.Item(9, 2).Value = 'BGC0001315'
$ws.Cells.Item(9, 3).Value = 'None'
$ws.Cells.Item(9, 4).Value = 'None'
$ws.Cells.Item(9, 5).Value = 'Cucumis sativus'
$ws.Cells.Item(9, 6).Value = 0
$ws.Cells.Item(9, 7).Value = 35297
$ws.Cells.Item(9, 8).Value = 'None'
$ws.Cells.Item(9, 9).Value = 6
$ws.Cells.Item(9, 10).Value = 'Terpene'
$ws.Cells.Item(9, 11).Value = 'Cucurbitacin biosynthetic gene cluster'
$ws.Cells.Item(9, 12).Value = 'cucurbitacin'

# Row 10: BGC0001316
$ws.Cells.Item(10, 1).Value = 8
$ws.Cells.Item(10, 2).Value = 'BGC0001316'
$ws.Cells.Item(10, 3).Value = 'None'
$ws.Cells.Item(10, 4).Value = 'None'
$ws.Cells.Item(10, 5).Value = 'Lotus japonicus'
$ws.Cells.Item(10, 6).Value = 0
$ws.Cells.Item(10, 7).Value = 409731
$ws.Cells.Item(10, 8).Value = 'None'
$ws.Cells.Item(10, 9).Value = 24
$ws.Cells.Item(10, 10).Value = 'Other'
$ws.Cells.Item(10, 11).Value = 'Linamarin / Lotaustralin biosynthetic gene cluster'
$ws.Cells.Item(10, 12).Value = 'linamarin,'

# Row 11: BGC0001317
$ws.Cells.Item(11, 1).Value = 9
$ws.Cells.Item(11, 2).Value = 'BGC0001317'
$ws.Cells.Item(11, 3).Value = 'None'
$ws.Cells.Item(11, 4).Value = 'None'
$ws.Cells.Item(11, 5).Value = 'Lotus japonicus'
$ws.Cells.Item(11, 6).Value = 0
$ws.Cells.Item(11, 7).Value = 277585
$ws.Cells.Item(11, 8).Value = 'None'
$ws.Cells.Item(11, 9).Value = 18
$ws.Cells.Item(11, 10).Value = 'Terpene'
$ws.Cells.Item(11, 11).Value = 'Lupeol biosynthetic gene cluster'
$ws.Cells.Item(11, 12).Value = 'lupeol'

# Row 12: BGC0001318
$ws.Cells.Item(12, 1).Value = 10
$ws.Cells.Item(12, 2).Value = 'BGC0001318'
$ws.Cells.Item(12, 3).Value = 'None'
$ws.Cells.Item(12, 4).Value = 'None'
$ws.Cells.Item(12, 5).Value = 'Manihot esculenta'
$ws.Cells.Item(12, 6).Value = 0
$ws.Cells.Item(12, 7).Value = 80994
$ws.Cells.Item(12, 8).Value = 'None'
$ws.Cells.Item(12, 9).Value = 13
$ws.Cells.Item(12, 10).Value = 'Other'
$ws.Cells.Item(12, 11).Value = 'Linamarin / Lotaustralin biosynthetic gene cluster'
$ws.Cells.Item(12, 12).Value = 'linamarin,'

# Row 13: BGC0001324
$ws.Cells.Item(13, 1).Value = 11
$ws.Cells.Item(13, 2).Value = 'BGC0001324'
$ws.Cells.Item(13, 3).Value = 'None'
$ws.Cells.Item(13, 4).Value = 'None'
$ws.Cells.Item(13, 5).Value = 'Solanum pimpinellifolium'
$ws.Cells.Item(13, 6).Value = 0
$ws.Cells.Item(13, 7).Value = 107124
$ws.Cells.Item(13, 8).Value = 'None'
$ws.Cells.Item(13, 9).Value = 0
$ws.Cells.Item(13, 10).Value = 'Terpene'
$ws.Cells.Item(13, 11).Value = 'Solanum pimpinellifolium isolate LA1589 chromosome 8 terpene biosynthesis gene locus, partial sequence'
$ws.Cells.Item(13, 12).Value = 'monoterpenes-diterpenes'

# Row 14: BGC0001325
$ws.Cells.Item(14, 1).Value = 12
$ws.Cells.Item(14, 2).Value = 'BGC0001325'
$ws.Cells.Item(14, 3).Value = 'None'
$ws.Cells.Item(14, 4).Value = 'None'
$ws.Cells.Item(14, 5).Value = 'Papaver somniferum'
$ws.Cells.Item(14, 6).Value = 0
$ws.Cells.Item(14, 7).Value = 401328
$ws.Cells.Item(14, 8).Value = 'None'
$ws.Cells.Item(14, 9).Value = 10
$ws.Cells.Item(14, 10).Value = 'Alkaloid'
$ws.Cells.Item(14, 11).Value = 'Noscapine biosynthetic gene cluster'
$ws.Cells.Item(14, 12).Value = 'noscapine'

# Row 15: BGC0001756
$ws.Cells.Item(15, 1).Value = 13
$ws.Cells.Item(15, 2).Value = 'BGC0001756'
$ws.Cells.Item(15, 3).Value = 'None'
$ws.Cells.Item(15, 4).Value = 'None'
$ws.Cells.Item(15, 5).Value = 'Arabidopsis thaliana'
$ws.Cells.Item(15, 6).Value = 4863612
$ws.Cells.Item(15, 7).Value = 4887487
$ws.Cells.Item(15, 8).Value = 'None'
$ws.Cells.Item(15, 9).Value = 10
$ws.Cells.Item(15, 10).Value = 'Terpene'
$ws.Cells.Item(15, 11).Value = 'Arabidopsis thaliana chromosome 3, partial sequence'
$ws.Cells.Item(15, 12).Value = 'arathanatriene,'

# Row 16: BGC0001799
$ws.Cells.Item(16, 1).Value = 14
$ws.Cells.Item(16, 2).Value = 'BGC0001799'
$ws.Cells.Item(16, 3).Value = 'None'
$ws.Cells.Item(16, 4).Value = 'None'
$ws.Cells.Item(16, 5).Value = 'Papaver somniferum'
$ws.Cells.Item(16, 6).Value = 78651
$ws.Cells.Item(16, 7).Value = 1013747
$ws.Cells.Item(16, 8).Value = 'None'
$ws.Cells.Item(16, 9).Value = 12
$ws.Cells.Item(16, 10).Value = 'Alkaloid'
$ws.Cells.Item(16, 11).Value = 'UNVERIFIED: Papaver somniferum (S)-reticuline epimerase-like (REPI1), REPI2, salutaridine synthase (SalSyn1), O-methyltransferase-1, SalSyn2, O-methyltransferase-2, salutaridinol 7-O-acetyltransferase (SalAT2), salutaridine reductase-like (SalR2), and thebaine synthase 2 (THS2) genes, partial sequence; thebaine synthase 1-like (THS1) gene, complete sequence; and SalR1 and SalAT1 genes, partial sequence'
$ws.Cells.Item(16, 12).Value = 'thebaine'

# Row 17: BGC0001997
$ws.Cells.Item(17, 1).Value = 15
$ws.Cells.Item(17, 2).Value = 'BGC0001997'
$ws.Cells.Item(17, 3).Value = 'None'
$ws.Cells.Item(17, 4).Value = 'None'
$ws.Cells.Item(17, 5).Value = 'Nicotiana tabacum'
$ws.Cells.Item(17, 6).Value = 469295
$ws.Cells.Item(17, 7).Value = 543944
$ws.Cells.Item(17, 8).Value = 'None'
$ws.Cells.Item(17, 9).Value = 0
$ws.Cells.Item(17, 10).Value = 'Terpene'
$ws.Cells.Item(17, 11).Value = 'Nicotiana tabacum cultivar K326 Nitab4.5_0001461, whole genome shotgun sequence'
$ws.Cells.Item(17, 12).Value = 'capsidiol'

# Row 18: BGC0002389
$ws.Cells.Item(18, 1).Value = 16
$ws.Cells.Item(18, 2).Value = 'BGC0002389'
$ws.Cells.Item(18, 3).Value = 'None'
$ws.Cells.Item(18, 4).Value = 'None'
$ws.Cells.Item(18, 5).Value = 'Zea mays'
$ws.Cells.Item(18, 6).Value = 56565232
$ws.Cells.Item(18, 7).Value = 56847613
$ws.Cells.Item(18, 8).Value = 'None'
$ws.Cells.Item(18, 9).Value = 4
$ws.Cells.Item(18, 10).Value = 'Terpene'
$ws.Cells.Item(18, 11).Value = 'Zea mays cultivar B73 chromosome 10, Zm-B73-REFERENCE-NAM-5.0, whole genome shotgun sequence'
$ws.Cells.Item(18, 12).Value = 'β-bisabolene,'

# Row 19: BGC0002390
$ws.Cells.Item(19, 1).Value = 17
$ws.Cells.Item(19, 2).Value = 'BGC0002390'
$ws.Cells.Item(19, 3).Value = 'None'
$ws.Cells.Item(19, 4).Value = 'None'
$ws.Cells.Item(19, 5).Value = 'Zea mays'
$ws.Cells.Item(19, 6).Value = 33367131
$ws.Cells.Item(19, 7).Value = 34010443
$ws.Cells.Item(19, 8).Value = 'None'
$ws.Cells.Item(19, 9).Value = 14
$ws.Cells.Item(19, 10).Value = 'Terpene'
$ws.Cells.Item(19, 11).Value = 'Zea mays cultivar B73 chromosome 5, Zm-B73-REFERENCE-NAM-5.0, whole genome shotgun sequence'
$ws.Cells.Item(19, 12).Value = 'zealexin'

# Row 20: BGC0002391
$ws.Cells.Item(20, 1).Value = 18
$ws.Cells.Item(20, 2).Value = 'BGC0002391'
$ws.Cells.Item(20, 3).Value = 'None'
$ws.Cells.Item(20, 4).Value = 'None'
$ws.Cells.Item(20, 5).Value = 'Zea mays'
$ws.Cells.Item(20, 6).Value = 285583113
$ws.Cells.Item(20, 7).Value = 285654925
$ws.Cells.Item(20, 8).Value = 'None'
$ws.Cells.Item(20, 9).Value = 2
$ws.Cells.Item(20, 10).Value = 'Terpene'
$ws.Cells.Item(20, 11).Value = 'Zea mays cultivar B73 chromosome 1, Zm-B73-REFERENCE-NAM-5.0, whole genome shotgun sequence'
$ws.Cells.Item(20, 12).Value = 'zealexin'

# Row 21: BGC0002392
$ws.Cells.Item(21, 1).Value = 19
$ws.Cells.Item(21, 2).Value = 'BGC0002392'
$ws.Cells.Item(21, 3).Value = 'None'
$ws.Cells.Item(21, 4).Value = 'None'
$ws.Cells.Item(21, 5).Value = 'Oryza sativa Japonica Group'
$ws.Cells.Item(21, 6).Value = 6494479
$ws.Cells.Item(21, 7).Value = 6634480
$ws.Cells.Item(21, 8).Value = 'None'
$ws.Cells.Item(21, 9).Value = 16
$ws.Cells.Item(21, 10).Value = 'Terpene'
$ws.Cells.Item(21, 11).Value = 'Oryza sativa Japonica Group DNA, chromosome 7, cultivar: Nipponbare, complete sequence'
$ws.Cells.Item(21, 12).Value = '5,10-diketo-casbene'

# Row 22: BGC0002393
$ws.Cells.Item(22, 1).Value = 20
$ws.Cells.Item(22, 2).Value = 'BGC0002393'
$ws.Cells.Item(22, 3).Value = 'None'
$ws.Cells.Item(22, 4).Value = 'None'
$ws.Cells.Item(22, 5).Value = 'Ricinus communis'
$ws.Cells.Item(22, 6).Value = 264511
$ws.Cells.Item(22, 7).Value = 337591
$ws.Cells.Item(22, 8).Value = 'None'
$ws.Cells.Item(22, 9).Value = 17
$ws.Cells.Item(22, 10).Value = 'Terpene'
$ws.Cells.Item(22, 11).Value = 'Ricinus communis genomic scaffold scf_1106159296192, whole genome shotgun sequence'
$ws.Cells.Item(22, 12).Value = '5a-hydroxy-casbene,'

# Row 23: BGC0002394
$ws.Cells.Item(23, 1).Value = 21
$ws.Cells.Item(23, 2).Value = 'BGC0002394'
$ws.Cells.Item(23, 3).Value = 'None'
$ws.Cells.Item(23, 4).Value = 'None'
$ws.Cells.Item(23, 5).Value = 'Taxus chinensis'
$ws.Cells.Item(23, 6).Value = 55305454
$ws.Cells.Item(23, 7).Value = 55566904
$ws.Cells.Item(23, 8).Value = 'None'
$ws.Cells.Item(23, 9).Value = 6
$ws.Cells.Item(23, 10).Value = 'Terpene'
$ws.Cells.Item(23, 11).Value = 'Taxus chinensis isolate Ta-2019 chromosome 9, whole genome shotgun sequence'
$ws.Cells.Item(23, 12).Value = 'taxa-4(20),11(12)-diene,'

# Row 24: BGC0002395
$ws.Cells.Item(24, 1).Value = 22
$ws.Cells.Item(24, 2).Value = 'BGC0002395'
$ws.Cells.Item(24, 3).Value = 'None'
$ws.Cells.Item(24, 4).Value = 'None'
$ws.Cells.Item(24, 5).Value = 'Hordeum vulgare subsp. vulgare'
$ws.Cells.Item(24, 6).Value = 9581241
$ws.Cells.Item(24, 7).Value = 10180879
$ws.Cells.Item(24, 8).Value = 'None'
$ws.Cells.Item(24, 9).Value = 19
$ws.Cells.Item(24, 10).Value = 'Terpene'
$ws.Cells.Item(24, 11).Value = 'Hordeum vulgare subsp. vulgare chromosome 2H, MorexV3_pseudomolecules_assembly, whole genome shotgun sequence'
$ws.Cells.Item(24, 12).Value = '11-hydroxy-hordetriene,'

# Row 25: BGC0002404
$ws.Cells.Item(25, 1).Value = 23
$ws.Cells.Item(25, 2).Value = 'BGC0002404'
$ws.Cells.Item(25, 3).Value = 'None'
$ws.Cells.Item(25, 4).Value = 'None'
$ws.Cells.Item(25, 5).Value = 'Solanum lycopersicum'
$ws.Cells.Item(25, 6).Value = 68007492
$ws.Cells.Item(25, 7).Value = 68031028
$ws.Cells.Item(25, 8).Value = 'None'
$ws.Cells.Item(25, 9).Value = 3
$ws.Cells.Item(25, 10).Value = 'Other'
$ws.Cells.Item(25, 11).Value = 'Solanum lycopersicum cultivar Heinz 1706 chromosome 12, SL3.0, whole genome shotgun sequence'
$ws.Cells.Item(25, 12).Value = 'falcarindiol'

# Row 26: BGC0002405
$ws.Cells.Item(26, 1).Value = 24
$ws.Cells.Item(26, 2).Value = 'BGC0002405'
$ws.Cells.Item(26, 3).Value = 'None'
$ws.Cells.Item(26, 4).Value = 'None'
$ws.Cells.Item(26, 5).Value = 'Solanum lycopersicum'
$ws.Cells.Item(26, 6).Value = 57739269
$ws.Cells.Item(26, 7).Value = 57782130
$ws.Cells.Item(26, 8).Value = 'None'
$ws.Cells.Item(26, 9).Value = 5
$ws.Cells.Item(26, 10).Value = 'Saccharide'
$ws.Cells.Item(26, 11).Value = 'Solanum lycopersicum cultivar Heinz 1706 chromosome 7, SL3.0, whole genome shotgun sequence'
$ws.Cells.Item(26, 12).Value = 'mid-chain'

# Row 27: BGC0002406
$ws.Cells.Item(27, 1).Value = 25
$ws.Cells.Item(27, 2).Value = 'BGC0002406'
$ws.Cells.Item(27, 3).Value = 'None'
$ws.Cells.Item(27, 4).Value = 'None'
$ws.Cells.Item(27, 5).Value = 'Oryza sativa Japonica Group'
$ws.Cells.Item(27, 6).Value = 12067615
$ws.Cells.Item(27, 7).Value = 12262361
$ws.Cells.Item(27, 8).Value = 'None'
$ws.Cells.Item(27, 9).Value = 17
$ws.Cells.Item(27, 10).Value = 'Other'
$ws.Cells.Item(27, 11).Value = 'Oryza sativa Japonica Group DNA, chromosome 10, cultivar: Nipponbare, complete sequence'
$ws.Cells.Item(27, 12).Value = 'hydroxycinnamoyltyramine'

# Row 28: BGC0002622
$ws.Cells.Item(28, 1).Value = 26
$ws.Cells.Item(28, 2).Value = 'BGC0002622'
$ws.Cells.Item(28, 3).Value = 'None'
$ws.Cells.Item(28, 4).Value = 'None'
$ws.Cells.Item(28, 5).Value = 'Oryza sativa Japonica Group'
$ws.Cells.Item(28, 6).Value = 21409012
$ws.Cells.Item(28, 7).Value = 21465464
$ws.Cells.Item(28, 8).Value = 'None'
$ws.Cells.Item(28, 9).Value = 3
$ws.Cells.Item(28, 10).Value = 'Other'
$ws.Cells.Item(28, 11).Value = 'Oryza sativa Japonica Group DNA, chromosome 9, cultivar: Nipponbare, complete sequence'
$ws.Cells.Item(28, 12).Value = 'hydroxycinnamoylputrescine'

# Row 29: BGC0002721
$ws.Cells.Item(29, 1).Value = 27
$ws.Cells.Item(29, 2).Value = 'BGC0002721'
$ws.Cells.Item(29, 3).Value = 'None'
$ws.Cells.Item(29, 4).Value = 'None'
$ws.Cells.Item(29, 5).Value = 'Hordeum vulgare subsp. vulgare'
$ws.Cells.Item(29, 6).Value = 17320184
$ws.Cells.Item(29, 7).Value = 17862889
$ws.Cells.Item(29, 8).Value = 'None'
$ws.Cells.Item(29, 9).Value = 5
$ws.Cells.Item(29, 10).Value = 'Saccharide'
$ws.Cells.Item(29, 11).Value = 'Hordeum vulgare subsp. vulgare chromosome 1H, MorexV3_pseudomolecules_assembly, whole genome shotgun sequence'
$ws.Cells.Item(29, 12).Value = 'α-hydroxynitrile'

# Row 30: BGC0002722
$ws.Cells.Item(30, 1).Value = 28
$ws.Cells.Item(30, 2).Value = 'BGC0002722'
$ws.Cells.Item(30, 3).Value = 'None'
$ws.Cells.Item(30, 4).Value = 'None'
$ws.Cells.Item(30, 5).Value = 'Saccharide'
$ws.Cells.Item(30, 6).Value = 1123004
$ws.Cells.Item(30, 7).Value = 1460949
$ws.Cells.Item(30, 8).Value = 'None'
$ws.Cells.Item(30, 9).Value = 22
$ws.Cells.Item(30, 10).Value = 'Alkaloid'
$ws.Cells.Item(30, 11).Value = 'Solanum tuberosum cultivar DM 1-3 516 R44 unplaced genomic scaffold, SolTub_3.0 scf00140, whole genome shotgun sequence'
$ws.Cells.Item(30, 12).Value = 'α-chaconine,'

# Row 31: BGC0002723
$ws.Cells.Item(31, 1).Value = 29
$ws.Cells.Item(31, 2).Value = 'BGC0002723'
$ws.Cells.Item(31, 3).Value = 'None'
$ws.Cells.Item(31, 4).Value = 'None'
$ws.Cells.Item(31, 5).Value = 'Jatropha curcas'
$ws.Cells.Item(31, 6).Value = 7856
$ws.Cells.Item(31, 7).Value = 18675
$ws.Cells.Item(31, 8).Value = 'None'
$ws.Cells.Item(31, 9).Value = 2
$ws.Cells.Item(31, 10).Value = 'Terpene'
$ws.Cells.Item(31, 11).Value = 'Jatropha curcas isolate ELGS0001-1231 unplaced genomic scaffold, RJC1_Hi-C scaffold_928, whole genome shotgun sequence'
$ws.Cells.Item(31, 12).Value = 'casbene'

# Row 32: BGC0002724
$ws.Cells.Item(32, 1).Value = 30
$ws.Cells.Item(32, 2).Value = 'BGC0002724'
$ws.Cells.Item(32, 3).Value = 'None'
$ws.Cells.Item(32, 4).Value = 'None'
$ws.Cells.Item(32, 5).Value = 'Jatropha curcas'
$ws.Cells.Item(32, 6).Value = 848628
$ws.Cells.Item(32, 7).Value = 912207
$ws.Cells.Item(32, 8).Value = 'None'
$ws.Cells.Item(32, 9).Value = 6
$ws.Cells.Item(32, 10).Value = 'Terpene'
$ws.Cells.Item(32, 11).Value = 'Jatropha curcas isolate ELGS0001-1231 unplaced genomic scaffold, RJC1_Hi-C scaffold_989, whole genome shotgun sequence'
$ws.Cells.Item(32, 12).Value = 'casbene'
